# Apply the 11/13/2019 5:38 pm update: add a "Model 7" results column (I)
# and refresh the accuracy numbers for Models 1-6 (columns C-H) across all
# sparsity rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column I: header + formatting ------------------------------------
# Column H's header cell (H1) carries the bold/centered/bordered header
# style; copy that formatting onto I1 before assigning its label so I1
# matches the other header cells.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "Model 7"

# --- Updated results for Models 1-6 (columns C-H), rows 2-11 --------------
$ws.Range("C2").Value = 0.6200000047683716
$ws.Range("D2").Value = 0.859499990940094
$ws.Range("E2").Value = 0.7789999842643738
$ws.Range("F2").Value = 0.7595000267028809
$ws.Range("G2").Value = 0.7925000190734863
$ws.Range("H2").Value = 0.8355000019073486

$ws.Range("C3").Value = 0.8999999761581421
$ws.Range("D3").Value = 0.9775000214576721
$ws.Range("E3").Value = 0.9704999923706055
$ws.Range("F3").Value = 0.9635000228881836
$ws.Range("G3").Value = 0.9695000052452087
$ws.Range("H3").Value = 0.9789999723434448

$ws.Range("C4").Value = 0.878000020980835
$ws.Range("D4").Value = 0.9950000047683716
$ws.Range("E4").Value = 0.9909999966621399
$ws.Range("F4").Value = 0.9835000038146973
$ws.Range("G4").Value = 0.9850000143051147
$ws.Range("H4").Value = 0.9955000281333923

$ws.Range("C5").Value = 0.8859999775886536
$ws.Range("D5").Value = 0.9965000152587891
$ws.Range("E5").Value = 0.987500011920929
$ws.Range("F5").Value = 0.9879999756813049
$ws.Range("G5").Value = 0.9890000224113464
$ws.Range("H5").Value = 0.9934999942779541

$ws.Range("C6").Value = 0.9925000071525574
$ws.Range("D6").Value = 0.5109999775886536
$ws.Range("E6").Value = 0.7269999980926514
$ws.Range("F6").Value = 0.815500020980835
$ws.Range("G6").Value = 0.5870000123977661
$ws.Range("H6").Value = 0.6104999780654907

$ws.Range("C7").Value = 0.9940000176429749
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.9994999766349792
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 1

$ws.Range("C8").Value = 0.9994999766349792
$ws.Range("D8").Value = 0.9980000257492065
$ws.Range("E8").Value = 0.9984999895095825
$ws.Range("F8").Value = 0.9994999766349792
$ws.Range("G8").Value = 0.9980000257492065
$ws.Range("H8").Value = 0.9944999814033508

$ws.Range("C9").Value = 0.9994999766349792
$ws.Range("D9").Value = 0.9990000128746033
$ws.Range("E9").Value = 0.9975000023841858
$ws.Range("F9").Value = 0.9994999766349792
$ws.Range("G9").Value = 0.9975000023841858
$ws.Range("H9").Value = 0.9975000023841858

$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 0.9994999766349792

$ws.Range("C11").Value = 0.9994999766349792
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 1

# --- New column I values, rows 2-11 ----------------------------------------
$ws.Range("I2").Value = 0.6769999861717224
$ws.Range("I3").Value = 0.9359999895095825
$ws.Range("I4").Value = 0.9620000123977661
$ws.Range("I5").Value = 0.9670000076293945
$ws.Range("I6").Value = 0.9980000257492065
$ws.Range("I7").Value = 0.9975000023841858
$ws.Range("I8").Value = 0.9994999766349792
$ws.Range("I9").Value = 0.9994999766349792
$ws.Range("I10").Value = 1
$ws.Range("I11").Value = 1
